$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/date range) ---
$ws.Range("A8").Value = "Volume 31   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/21/2024  Through  10/27/2024"

# --- Crime-data table updates (rows 14-28) ---
# Type-changing cells: set value (forcing text with leading quote where needed), then fix style by pasting format from a same-style donor cell
$ws.Range("C15").Value = "'0"
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("D16").Value = 4
$ws.Range("F16").Copy()
$ws.Range("D16").PasteSpecial(-4122)

$ws.Range("E16").Value = -25
$ws.Range("E17").Copy()
$ws.Range("E16").PasteSpecial(-4122)

$ws.Range("C20").Value = 1
$ws.Range("D20").Copy()
$ws.Range("C20").PasteSpecial(-4122)

$ws.Range("C22").Value = 1
$ws.Range("F22").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("D22").Value = 1
$ws.Range("F22").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").Value = 0
$ws.Range("H22").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("C27").Value = "'0"
$ws.Range("D27").Copy()
$ws.Range("C27").PasteSpecial(-4122)

# Plain numeric value updates (no type/style change)
$ws.Range("M14").Value = -66.666666666666
$ws.Range("C16").Value = 3
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 108
$ws.Range("J16").Value = 102
$ws.Range("K16").Value = 5.882352941176
$ws.Range("L16").Value = -0.917431192660
$ws.Range("M16").Value = -45.728643216080
$ws.Range("N16").Value = -81.25
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -28.571428571428
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 207
$ws.Range("J17").Value = 204
$ws.Range("K17").Value = 1.470588235294
$ws.Range("L17").Value = 7.8125
$ws.Range("M17").Value = 38.926174496644
$ws.Range("N17").Value = -64.186851211072
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 8
$ws.Range("H18").Value = 33.333333333333
$ws.Range("I18").Value = 75
$ws.Range("J18").Value = 56
$ws.Range("K18").Value = 33.928571428571
$ws.Range("L18").Value = 29.310344827586
$ws.Range("M18").Value = -19.354838709677
$ws.Range("N18").Value = -88.335925349922
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -38.461538461538
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = 14.814814814814
$ws.Range("I19").Value = 285
$ws.Range("J19").Value = 292
$ws.Range("K19").Value = -2.397260273972
$ws.Range("L19").Value = -3.061224489795
$ws.Range("M19").Value = 133.606557377049
$ws.Range("N19").Value = 1.423487544483
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = -46.153846153846
$ws.Range("I20").Value = 64
$ws.Range("J20").Value = 97
$ws.Range("K20").Value = -34.020618556701
$ws.Range("L20").Value = -4.477611940298
$ws.Range("M20").Value = 48.837209302325
$ws.Range("N20").Value = -75.849056603773
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -25
$ws.Range("G21").Value = 76
$ws.Range("H21").Value = 3.947368421052
$ws.Range("I21").Value = 749
$ws.Range("J21").Value = 759
$ws.Range("K21").Value = -1.317523056653
$ws.Range("L21").Value = 1.216216216216
$ws.Range("M21").Value = 17.952755905511
$ws.Range("N21").Value = -69.240246406570
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 21
$ws.Range("J22").Value = 28
$ws.Range("K22").Value = -25
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 90.909090909090
$ws.Range("C24").Value = 10
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = -37.5
$ws.Range("F24").Value = 71
$ws.Range("G24").Value = 52
$ws.Range("H24").Value = 36.538461538461
$ws.Range("I24").Value = 636
$ws.Range("J24").Value = 655
$ws.Range("K24").Value = -2.900763358778
$ws.Range("L24").Value = -36.336336336336
$ws.Range("M24").Value = 92.145015105740
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 6
$ws.Range("H25").Value = 350
$ws.Range("I25").Value = 186
$ws.Range("J25").Value = 159
$ws.Range("K25").Value = 16.981132075471
$ws.Range("L25").Value = -64.023210831721
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 140
$ws.Range("F26").Value = 36
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = 44
$ws.Range("I26").Value = 304
$ws.Range("J26").Value = 263
$ws.Range("K26").Value = 15.589353612167
$ws.Range("L26").Value = 3.050847457627
$ws.Range("M26").Value = -27.619047619047
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 2
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 24
$ws.Range("J28").Value = 27
$ws.Range("K28").Value = -11.111111111111
$ws.Range("L28").Value = -51.020408163265

$ws.Application.CutCopyMode = $false
Write-Host "edit applied"
